# Saldo.xlsx update
# - Row 4 (account 004001621, DANIELA): Saldo updated from 100484.68 to 165512.75
# - Remove row for account 003301389 (EDMUNDO, 70000)
# - Remove row for account 004460491 (PEDRO, 4988.96)
# - Remove row for account 004205505 (SURAMA, 1228.26)
# - Remove row for account 004870019 (MARIA, 768.38)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DANIELA's balance (row 4, column C)
$ws.Cells.Item(4, 3).Value = 165512.75

# Delete the obsolete rows. Delete from the bottom up so earlier row
# numbers stay valid as each row is removed.
$ws.Rows.Item(40).EntireRow.Delete()
$ws.Rows.Item(14).EntireRow.Delete()
$ws.Rows.Item(11).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()
